# POSt Of reconciliation Steps
# Update the bank-statement reconciliation sheet: refresh the From/To dates,
# refresh the existing two transactions, and post the same two transactions
# again for the new statement period (rows 8-9), extending the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style touch-ups -------------------------------------------------
# Row 5 (column headers) switch to the bold "Consolas" label style already
# used by the "Bank statement:" caption in A1.
$ws.Range("A1").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)   # xlPasteFormats

# D5/G5 (Details / STATUS headers) share one consolidated centered style.
$ws.Range("G5").Copy()
$ws.Range("D5").PasteSpecial(-4122)

# C8 takes the same "centered number" style used by A6/A7/C6/C7.
$ws.Range("A6").Copy()
$ws.Range("C8").PasteSpecial(-4122)

# G8 drops the stray custom number format and matches G6/G7's plain style.
$ws.Range("G6").Copy()
$ws.Range("G8").PasteSpecial(-4122)

# Row 9 is a brand-new row: give it the same per-column look as row 7
# (transaction ID / date / reference / details / debit / credit / status).
$ws.Range("A7:G7").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Header dates ------------------------------------------------------
$ws.Range("C2").Value2 = 45075
$ws.Range("C3").Value2 = 45077

# --- Row 6: first transaction refreshed --------------------------------
$ws.Range("A6").Value2 = 25545
$ws.Range("B6").Value2 = 45075
$ws.Range("C6").Value2 = 2034
$ws.Range("F6").Value2 = 240

# --- Row 7: second transaction refreshed -------------------------------
$ws.Range("A7").Value2 = 25546
$ws.Range("B7").Value2 = 45075
$ws.Range("E7").Value2 = 20

# --- Row 8: first transaction posted again for the new period ----------
$ws.Range("A8").Value2 = 25545
$ws.Range("B8").Value2 = 45077
$ws.Range("C8").Value2 = 2034
$ws.Range("D8").Value2 = "Incoming Payment"
$ws.Range("F8").Value2 = 12000

# --- Row 9 (new): second transaction posted again for the new period ---
$ws.Range("A9").Value2 = 25546
$ws.Range("B9").Value2 = 45077
$ws.Range("D9").Value2 = "Bank Charge"
$ws.Range("E9").Value2 = 250
